# Update forecast workbook with corrected forecast output.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Forecast Comparison": insert a new "Week_Start_Date" column
# after "Week" (new column B), shifting everything right by one column,
# renumber week labels (W01 -> W1, etc.), update MyForecast values, and
# mark the (now) last column as boolean is_holiday_week.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Insert new column B ("Week_Start_Date"); everything from old column B
# onward shifts one column to the right (B->C, C->D, ... H->I, I->J).
$ws1.Columns.Item(2).Insert()

# Header row
$ws1.Cells.Item(1, 2).Value = "Week_Start_Date"

# Per-row data: Week label, Week_Start_Date, MyForecast (new column D)
$weekData = @(
    @{ Row = 2;  Week = "W1";  Start = "2025-01-05"; MyForecast = 67 },
    @{ Row = 3;  Week = "W2";  Start = "2025-01-12"; MyForecast = 64 },
    @{ Row = 4;  Week = "W3";  Start = "2025-01-19"; MyForecast = 67 },
    @{ Row = 5;  Week = "W4";  Start = "2025-01-26"; MyForecast = 63 },
    @{ Row = 6;  Week = "W5";  Start = "2025-02-02"; MyForecast = 61 },
    @{ Row = 7;  Week = "W6";  Start = "2025-02-09"; MyForecast = 59 },
    @{ Row = 8;  Week = "W7";  Start = "2025-02-16"; MyForecast = 64 },
    @{ Row = 9;  Week = "W8";  Start = "2025-02-23"; MyForecast = 75 },
    @{ Row = 10; Week = "W9";  Start = "2025-03-02"; MyForecast = 55 },
    @{ Row = 11; Week = "W10"; Start = "2025-03-09"; MyForecast = 61 },
    @{ Row = 12; Week = "W11"; Start = "2025-03-16"; MyForecast = 58 },
    @{ Row = 13; Week = "W12"; Start = "2025-03-23"; MyForecast = 59 },
    @{ Row = 14; Week = "W13"; Start = "2025-03-30"; MyForecast = 69 },
    @{ Row = 15; Week = "W14"; Start = "2025-04-06"; MyForecast = 65 },
    @{ Row = 16; Week = "W15"; Start = "2025-04-13"; MyForecast = 61 },
    @{ Row = 17; Week = "W16"; Start = "2025-04-20"; MyForecast = 49 }
)

foreach ($item in $weekData) {
    $r = $item.Row

    # A: Week label (drop leading zero, e.g. W01 -> W1)
    $ws1.Cells.Item($r, 1).Value = $item.Week

    # B: Week_Start_Date, stored as plain text (not an Excel date)
    $cellB = $ws1.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $item.Start

    # D: MyForecast (corrected values)
    $ws1.Cells.Item($r, 4).Value = $item.MyForecast

    # J: is_holiday_week, now a boolean cell
    $ws1.Cells.Item($r, 10).Value = $false
}

# ---------------------------------------------------------------------
# Sheet "Summary": update recalculated forecast totals
# All Value column cells are plain text, so force text number format
# before assigning so Excel doesn't reinterpret them as numbers/dates.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$summaryUpdates = @(
    @{ Row = 9;  Text = "996" },
    @{ Row = 10; Text = "519" },
    @{ Row = 11; Text = "261" },
    @{ Row = 12; Text = "75" },
    @{ Row = 13; Text = "2025-02-23" }
)

foreach ($item in $summaryUpdates) {
    $cell = $ws2.Cells.Item($item.Row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Text
}
